# "Added Priority to Login&Queue"
# - Fix the "Please fill out this field" validation message on the Login
#   sheet so it ends with a period (matches the "EmptyValues" message used
#   on the Register sheet), which also drops the now-unused duplicate
#   shared string.
# - Make "Login" the prioritized / active sheet instead of "Queue", and
#   move the selection on Login to D4.

$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("Login")
$login.Range("D2").Value = "Please fill out this field."
$login.Range("D3").Value = "Please fill out this field."
$login.Range("D4").Value = "Please fill out this field."

$login.Activate()
$login.Range("D4").Select()
